$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 202, shifting existing rows 202:274 down to 203:275
$ws.Rows.Item(202).EntireRow.Insert()

# Populate the newly inserted row 202 with the new weekly record
$ws.Cells.Item(202, 1).Value = 9
$ws.Cells.Item(202, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(202, 3).Value = "Metropolitana"
$ws.Cells.Item(202, 4).Value = 44559
$ws.Cells.Item(202, 5).Value = 13
$ws.Cells.Item(202, 6).Value = 100112044
$ws.Cells.Item(202, 7).Value = "Perejil"
$ws.Cells.Item(202, 8).Value = "Sin especificar"
$ws.Cells.Item(202, 9).Value = "Primera"
$ws.Cells.Item(202, 10).Value = 79
$ws.Cells.Item(202, 11).Value = 10000
$ws.Cells.Item(202, 12).Value = 12000
$ws.Cells.Item(202, 13).Value = 10987
$ws.Cells.Item(202, 14).Value = "$/docena de atados"
$ws.Cells.Item(202, 15).Value = "Región Metropolitana"
$ws.Cells.Item(202, 16).Value = 3662
$ws.Cells.Item(202, 17).Value = 3
$ws.Cells.Item(202, 18).Value = "Hortaliza"

# Match the date number format used by the rest of column D
$ws.Cells.Item(202, 4).NumberFormat = $ws.Cells.Item(203, 4).NumberFormat
